$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.198371767997742
$ws.Range("B1").Value = 2.126335620880127
$ws.Range("C1").Value = 4.44853401184082
$ws.Range("D1").Value = 2.83938193321228
$ws.Range("E1").Value = 1.206780910491943
